# Localization status report refresh: items that were "Ready for handoff"
# have now moved to "In Translation". Update the status text everywhere it
# appears (Overview rollup columns + each per-locale sheet), then re-fit the
# status column width to the new (shorter) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Overview sheet: per-locale status rollup columns (E = zh-cn, F = de-de)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Per-locale sheets: "Status" column (column C)
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# The status column is narrower now that "In Translation" (14 chars) is
# shorter than "Ready for handoff" (17 chars) - re-fit the affected columns.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
